# Regenerate save_data to use K instead of Strike#, recalculated values for column G (K)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 2
    4  = 2
    5  = 2
    6  = 1
    7  = 4
    8  = 3
    9  = 0
    10 = 2
    11 = 1
    12 = 6
    13 = 3
    14 = 5
    15 = 3
    16 = 4
    17 = 6
    18 = 1
    19 = 2
    20 = 7
    21 = 1
    22 = 2
    23 = 3
    24 = 6
    25 = 3
    26 = 3
    27 = 7
    28 = 4
    29 = 6
    30 = 2
    31 = 5
    32 = 4
    33 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
